$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (E). Excel's entire-column delete shifts
# everything to the right of it (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) one column to the
# left, which matches the diff exactly.
$ws.Range("E:E").EntireColumn.Delete()
